# corrected issue related to importing BDF events as edftype or type
# idiosyncratically. Now gives visual feedback related to shifted event
# codes in bins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "EventsToAdjust" (F2) used to list the bin codes before the fix
# (112; 114; 122; 124). The fix adds the missing 132/134 codes so the
# shifted event codes show up in the bin listing.
$ws.Range("F2").Value = "112; 114; 122; 124; 132; 134"

# Reflect where the user was looking when they made/saved the change:
# scrolled right so column D is the left-most visible column, with the
# cursor sitting on the corrected cell (F2).
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollColumn = 4
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
$ws.Range("F2").Select()
